# Apply updated crypto price/volume data per GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.095.61"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "2.468.89"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'519.73"
$ws.Range("E5").Value = "  -3.00%  "
$ws.Range("D6").Value = "'134.17"
$ws.Range("E6").Value = "  -1.70%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("D9").Value = "2.479.10"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").Value = "'0.0987"
$ws.Range("E10").Value = "  -2.84%  "
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("D12").Value = "'5.27"
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").Value = "'0.340"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "2.907.00"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "58.056.49"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "'22.44"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "2.471.90"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "'10.71"
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("D20").Value = "'321.56"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("E21").Value = "  -1.81%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("D24").Value = "'64.41"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'0.159"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("D30").Value = "'168.20"
$ws.Range("E30").Value = "  -1.49%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.69"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'6.29"
$ws.Range("E32").Value = "  -5.14%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.18"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'18.08"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("D37").Value = "'1.32"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").Value = "'4.00"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("E42").Value = "  -3.05%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'5.09"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'273.67"
$ws.Range("E44").Value = "  -3.69%  "
$ws.Range("D45").Value = "'0.590"
$ws.Range("E45").Value = "  -2.76%  "
$ws.Range("D46").Value = "'124.44"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "'0.0911"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").Value = "'0.0490"
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("D49").Value = "'0.0213"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "'17.04"
$ws.Range("E50").Value = "  -2.00%  "
$ws.Range("D51").Value = "1.730.99"
